# Auto-generated edit script: updates computed market-price columns (H:N)
# for several leves across multiple sheets, per the scheduled runner refresh.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(55, 8).Value = 600
$ws.Cells.Item(55, 9).Value = 736.6667
$ws.Cells.Item(55, 11).Value = 736.6667
$ws.Cells.Item(55, 13).Value = -522.6667
$ws.Cells.Item(98, 8).Value = 1837
$ws.Cells.Item(98, 9).Value = 2208.889
$ws.Cells.Item(98, 10).Value = 1000.25
$ws.Cells.Item(98, 11).Value = 2208.889
$ws.Cells.Item(98, 12).Value = 1000.25
$ws.Cells.Item(98, 13).Value = -710.8890000000001
$ws.Cells.Item(98, 14).Value = -3996.25
$ws.Cells.Item(107, 8).Value = 8623262
$ws.Cells.Item(107, 9).Value = 10417494
$ws.Cells.Item(107, 10).Value = 10950
$ws.Cells.Item(107, 11).Value = 10417494
$ws.Cells.Item(107, 12).Value = 10950
$ws.Cells.Item(107, 13).Value = -10415574
$ws.Cells.Item(107, 14).Value = -14790
$ws.Cells.Item(112, 8).Value = 10870479
$ws.Cells.Item(112, 10).Value = 13514477
$ws.Cells.Item(112, 12).Value = 40543431
$ws.Cells.Item(112, 14).Value = -40545647
$ws.Cells.Item(122, 8).Value = 1837
$ws.Cells.Item(122, 9).Value = 2208.889
$ws.Cells.Item(122, 10).Value = 1000.25
$ws.Cells.Item(122, 11).Value = 6626.667
$ws.Cells.Item(122, 12).Value = 3000.75
$ws.Cells.Item(122, 13).Value = -4176.667
$ws.Cells.Item(122, 14).Value = -7900.75
$ws.Cells.Item(137, 8).Value = 1609.3024
$ws.Cells.Item(137, 9).Value = 1181.7778
$ws.Cells.Item(137, 10).Value = 1917.12
$ws.Cells.Item(137, 11).Value = 3545.3334
$ws.Cells.Item(137, 12).Value = 5751.36
$ws.Cells.Item(137, 13).Value = -995.3334000000004
$ws.Cells.Item(137, 14).Value = -10851.36

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 660.125
$ws.Cells.Item(110, 9).Value = 660.125
$ws.Cells.Item(110, 11).Value = 660.125
$ws.Cells.Item(110, 13).Value = 1384.875

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(75, 8).Value = 20000
$ws.Cells.Item(75, 9).Value = 20000
$ws.Cells.Item(75, 11).Value = 20000
$ws.Cells.Item(75, 13).Value = -19064
$ws.Cells.Item(78, 8).Value = 20000
$ws.Cells.Item(78, 9).Value = 20000
$ws.Cells.Item(78, 11).Value = 60000
$ws.Cells.Item(78, 13).Value = -55320
$ws.Cells.Item(107, 8).Value = 860.875
$ws.Cells.Item(107, 9).Value = 788.5833
$ws.Cells.Item(107, 10).Value = 1077.75
$ws.Cells.Item(107, 11).Value = 788.5833
$ws.Cells.Item(107, 12).Value = 1077.75
$ws.Cells.Item(107, 13).Value = 1131.4167
$ws.Cells.Item(107, 14).Value = -4917.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 9616127
$ws.Cells.Item(16, 9).Value = 12821201
$ws.Cells.Item(16, 11).Value = 12821201
$ws.Cells.Item(16, 13).Value = -12820914
$ws.Cells.Item(31, 8).Value = 4197.609
$ws.Cells.Item(31, 10).Value = 4829.1924
$ws.Cells.Item(31, 12).Value = 4829.1924
$ws.Cells.Item(31, 14).Value = -5419.1924
$ws.Cells.Item(34, 8).Value = 4197.609
$ws.Cells.Item(34, 10).Value = 4829.1924
$ws.Cells.Item(34, 12).Value = 4829.1924
$ws.Cells.Item(34, 14).Value = -5233.1924
$ws.Cells.Item(58, 8).Value = 1493.875
$ws.Cells.Item(58, 9).Value = 861.4286
$ws.Cells.Item(58, 10).Value = 2701.2727
$ws.Cells.Item(58, 11).Value = 861.4286
$ws.Cells.Item(58, 12).Value = 2701.2727
$ws.Cells.Item(58, 13).Value = -658.4286
$ws.Cells.Item(58, 14).Value = -3107.2727
$ws.Cells.Item(107, 8).Value = 552.1667
$ws.Cells.Item(107, 9).Value = 442.78946
$ws.Cells.Item(107, 10).Value = 674.41174
$ws.Cells.Item(107, 11).Value = 442.78946
$ws.Cells.Item(107, 12).Value = 674.41174
$ws.Cells.Item(107, 13).Value = 1477.21054
$ws.Cells.Item(107, 14).Value = -4514.41174
$ws.Cells.Item(113, 8).Value = 9616127
$ws.Cells.Item(113, 9).Value = 12821201
$ws.Cells.Item(113, 11).Value = 12821201
$ws.Cells.Item(113, 13).Value = -12819031
$ws.Cells.Item(136, 8).Value = 1493.875
$ws.Cells.Item(136, 9).Value = 861.4286
$ws.Cells.Item(136, 10).Value = 2701.2727
$ws.Cells.Item(136, 11).Value = 2584.2858
$ws.Cells.Item(136, 12).Value = 8103.8181
$ws.Cells.Item(136, 13).Value = -34.28579999999965
$ws.Cells.Item(136, 14).Value = -13203.8181

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 857.2646999999999
$ws.Cells.Item(107, 9).Value = 305
$ws.Cells.Item(107, 10).Value = 1104.0212
$ws.Cells.Item(107, 11).Value = 915
$ws.Cells.Item(107, 12).Value = 3312.063599999999
$ws.Cells.Item(107, 13).Value = 1005
$ws.Cells.Item(107, 14).Value = -7152.063599999999
$ws.Cells.Item(136, 8).Value = 16231.714
$ws.Cells.Item(136, 9).Value = 20624.4
$ws.Cells.Item(136, 10).Value = 5250
$ws.Cells.Item(136, 11).Value = 61873.2
$ws.Cells.Item(136, 12).Value = 15750
$ws.Cells.Item(136, 13).Value = -56773.2
$ws.Cells.Item(136, 14).Value = -25950
$ws.Cells.Item(137, 8).Value = 30318570
$ws.Cells.Item(137, 9).Value = 1517
$ws.Cells.Item(137, 10).Value = 38480852
$ws.Cells.Item(137, 11).Value = 4551
$ws.Cells.Item(137, 12).Value = 115442556
$ws.Cells.Item(137, 13).Value = 549
$ws.Cells.Item(137, 14).Value = -115452756
$ws.Cells.Item(140, 8).Value = 6608.36
$ws.Cells.Item(140, 9).Value = 6608.36
$ws.Cells.Item(140, 11).Value = 19825.08
$ws.Cells.Item(140, 13).Value = -14645.08
$ws.Cells.Item(141, 8).Value = 15846.263
$ws.Cells.Item(141, 9).Value = 12107.182
$ws.Cells.Item(141, 11).Value = 36321.546
$ws.Cells.Item(141, 13).Value = -31141.546

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 749.94116
$ws.Cells.Item(97, 9).Value = 768.1818
$ws.Cells.Item(97, 10).Value = 716.5
$ws.Cells.Item(97, 11).Value = 768.1818
$ws.Cells.Item(97, 12).Value = 716.5
$ws.Cells.Item(97, 13).Value = -272.1818
$ws.Cells.Item(97, 14).Value = -1708.5
$ws.Cells.Item(107, 8).Value = 2311.0527
$ws.Cells.Item(107, 9).Value = 10000
$ws.Cells.Item(107, 10).Value = 1883.8889
$ws.Cells.Item(107, 11).Value = 10000
$ws.Cells.Item(107, 12).Value = 1883.8889
$ws.Cells.Item(107, 13).Value = -8080
$ws.Cells.Item(107, 14).Value = -5723.8889
$ws.Cells.Item(109, 8).Value = 10285
$ws.Cells.Item(109, 10).Value = 10285
$ws.Cells.Item(109, 12).Value = 10285
$ws.Cells.Item(109, 14).Value = -12365
$ws.Cells.Item(133, 8).Value = 62296.668
$ws.Cells.Item(133, 10).Value = 62296.668
$ws.Cells.Item(133, 12).Value = 62296.668
$ws.Cells.Item(133, 14).Value = -72416.66800000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 5435109.5
$ws.Cells.Item(122, 9).Value = 7941670
$ws.Cells.Item(122, 10).Value = 1675268
$ws.Cells.Item(122, 11).Value = 23825010
$ws.Cells.Item(122, 12).Value = 5025804
$ws.Cells.Item(122, 13).Value = -23822560
$ws.Cells.Item(122, 14).Value = -5030704
$ws.Cells.Item(136, 8).Value = 4656.6665
$ws.Cells.Item(136, 9).Value = 4168.9766
$ws.Cells.Item(136, 10).Value = 5568.4346
$ws.Cells.Item(136, 11).Value = 12506.9298
$ws.Cells.Item(136, 12).Value = 16705.3038
$ws.Cells.Item(136, 13).Value = -9956.9298
$ws.Cells.Item(136, 14).Value = -21805.3038

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 58824016
$ws.Cells.Item(107, 10).Value = 676
$ws.Cells.Item(107, 12).Value = 2028
$ws.Cells.Item(107, 14).Value = -5868
$ws.Cells.Item(126, 8).Value = 1219.5
$ws.Cells.Item(126, 9).Value = 866.55554
$ws.Cells.Item(126, 10).Value = 1673.2858
$ws.Cells.Item(126, 11).Value = 2599.66662
$ws.Cells.Item(126, 12).Value = 5019.857400000001
$ws.Cells.Item(126, 13).Value = -129.66662
$ws.Cells.Item(126, 14).Value = -9959.857400000001
$ws.Cells.Item(136, 8).Value = 1235.8718
$ws.Cells.Item(136, 9).Value = 730.85187
$ws.Cells.Item(136, 10).Value = 2372.1667
$ws.Cells.Item(136, 11).Value = 2192.55561
$ws.Cells.Item(136, 12).Value = 7116.500100000001
$ws.Cells.Item(136, 13).Value = 357.4443900000001
$ws.Cells.Item(136, 14).Value = -12216.5001
